# Update F-column (热度/hits counters) values across sheets 展览, 演出, 全部类型
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsExhibit.Range("F3").Value = 256
$wsExhibit.Range("F4").Value = 261
$wsExhibit.Range("F5").Value = 1799
$wsExhibit.Range("F6").Value = 677
$wsExhibit.Range("F7").Value = 537
$wsExhibit.Range("F8").Value = 4931
$wsExhibit.Range("F9").Value = 59
$wsExhibit.Range("F12").Value = 993
$wsExhibit.Range("F14").Value = 1300
$wsExhibit.Range("F16").Value = 1888
$wsExhibit.Range("F17").Value = 3003
$wsExhibit.Range("F18").Value = 1844
$wsExhibit.Range("F20").Value = 50
$wsExhibit.Range("F22").Value = 72
$wsExhibit.Range("F23").Value = 654
$wsExhibit.Range("F25").Value = 314
$wsExhibit.Range("F27").Value = 3140
$wsExhibit.Range("F28").Value = 1044
$wsExhibit.Range("F29").Value = 2544
$wsExhibit.Range("F30").Value = 260
$wsExhibit.Range("F31").Value = 1388
$wsExhibit.Range("F32").Value = 3718
$wsExhibit.Range("F33").Value = 98
$wsExhibit.Range("F34").Value = 902
$wsExhibit.Range("F35").Value = 436
$wsExhibit.Range("F36").Value = 1163
$wsExhibit.Range("F37").Value = 7
$wsExhibit.Range("F38").Value = 945
$wsExhibit.Range("F39").Value = 1207
$wsExhibit.Range("F40").Value = 29
$wsExhibit.Range("F41").Value = 906
$wsExhibit.Range("F42").Value = 584
$wsExhibit.Range("F43").Value = 342
$wsExhibit.Range("F44").Value = 375
$wsExhibit.Range("F45").Value = 295
$wsExhibit.Range("F46").Value = 3508

$wsShow = $wb.Worksheets.Item(2)       # 演出
$wsShow.Range("F3").Value = 18
$wsShow.Range("F4").Value = 6

$wsAll = $wb.Worksheets.Item(4)        # 全部类型
$wsAll.Range("F3").Value = 256
$wsAll.Range("F4").Value = 261
$wsAll.Range("F6").Value = 1799
$wsAll.Range("F7").Value = 677
$wsAll.Range("F8").Value = 537
$wsAll.Range("F9").Value = 4931
$wsAll.Range("F10").Value = 59
$wsAll.Range("F11").Value = 18
$wsAll.Range("F14").Value = 1300
$wsAll.Range("F15").Value = 3003
$wsAll.Range("F17").Value = 1844
$wsAll.Range("F19").Value = 50
$wsAll.Range("F25").Value = 72
$wsAll.Range("F27").Value = 314
$wsAll.Range("F28").Value = 3140
$wsAll.Range("F30").Value = 1044
$wsAll.Range("F31").Value = 2544
$wsAll.Range("F32").Value = 1388
$wsAll.Range("F33").Value = 3718
$wsAll.Range("F34").Value = 98
$wsAll.Range("F35").Value = 902
$wsAll.Range("F36").Value = 1163
$wsAll.Range("F37").Value = 945
$wsAll.Range("F39").Value = 1207
$wsAll.Range("F40").Value = 29
$wsAll.Range("F41").Value = 906
$wsAll.Range("F42").Value = 584
$wsAll.Range("F43").Value = 375
$wsAll.Range("F47").Value = 295
$wsAll.Range("F48").Value = 3508

